$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at R for the new category "Entre núcleos de población
# ejidal o comunal y colonia agrícola que adopta dominio pleno" - this shifts
# every column from R onward one position to the right (R->S, S->T, ... AQ->AR).
$ws.Columns("R:R").Insert()

# Update the header text for column Q (clarified wording) and set the new
# column R header.
$ws.Range("Q1").Value = "Entre núcleos de población ejidal o comunal y sociedades o asociaciones"
$ws.Range("R1").Value = "Entre núcleos de población ejidal o comunal y colonia agrícola que adopta dominio pleno "

# New column R data for the two existing data rows.
$ws.Range("R2").Value = 0
$ws.Range("R3").Value = 0

# Row 2 (Distrito 34): period label casing fix.
$ws.Range("D2").Value = "Diciembre/2024"

# Row 2 (Distrito 34): updated figures.
$ws.Range("M2").Value = 105
$ws.Range("S2").Value = 1
$ws.Range("AB2").Value = 2
$ws.Range("AC2").Value = 4
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 103
$ws.Range("AP2").Value = 0
$ws.Range("AQ2").Value = 2
$ws.Range("AR2").Value = 0

# Row 3 (Distrito 34-A): updated figures.
$ws.Range("M3").Value = 105
$ws.Range("S3").Value = 1
$ws.Range("AB3").Value = 2
$ws.Range("AC3").Value = 4
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 103
$ws.Range("AP3").Value = 0
$ws.Range("AQ3").Value = 2
$ws.Range("AR3").Value = 0
